# Refresh the crypto price/volume snapshot (GitHub Actions scrape).
# Column D ("Price") holds plain text in the source sheet (values such as
# "46.895.50" or trailing-zero prices like "298.70" are not valid Excel
# numbers / would lose their formatting if auto-converted), so numeric-looking
# replacements are written with a leading apostrophe - Excel's standard
# "force text" entry convention - to keep them as text, exactly like the rest
# of the column. A doubled leading quote in a single-quoted PowerShell string
# (i.e. '''298.70') yields the literal text 'X, not a stray quote.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '46.895.50'
$ws.Range("E2").Value = '  +7.33%  '

# Row 3
$ws.Range("D3").Value = '2.312.83'
$ws.Range("E3").Value = '  +6.03%  '

# Row 4
$ws.Range("E4").Value = '  -0.44%  '

# Row 5
$ws.Range("D5").Value = '''298.70'
$ws.Range("E5").Value = '  +2.31%  '

# Row 6
$ws.Range("D6").Value = '''99.03'
$ws.Range("E6").Value = '  +14.52%  '

# Row 7
$ws.Range("D7").Value = '''0.571'
$ws.Range("E7").Value = '  +2.32%  '

# Row 8
$ws.Range("E8").Value = '  -0.43%  '

# Row 9
$ws.Range("D9").Value = '''0.529'
$ws.Range("E9").Value = '  +11.49%  '

# Row 10
$ws.Range("D10").Value = '''35.75'
$ws.Range("E10").Value = '  +12.63%  '

# Row 11
$ws.Range("E11").Value = '  +5.54%  '

# Row 12
$ws.Range("D12").Value = '''7.33'
$ws.Range("E12").Value = '  +10.33%  '

# Row 13
$ws.Range("E13").Value = '  +2.08%  '

# Row 14
$ws.Range("D14").Value = '2.663.55'
$ws.Range("E14").Value = '  +5.96%  '

# Row 15
$ws.Range("D15").Value = '2.308.95'
$ws.Range("E15").Value = '  +2.40%  '

# Row 16
$ws.Range("D16").Value = '''14.00'
$ws.Range("E16").Value = '  +9.48%  '

# Row 17
$ws.Range("D17").Value = '''0.815'
$ws.Range("E17").Value = '  +8.22%  '

# Row 18
$ws.Range("D18").Value = '46.815.72'
$ws.Range("E18").Value = '  +8.31%  '

# Row 19
$ws.Range("D19").Value = '''13.13'
$ws.Range("E19").Value = '  +24.41%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0940'
$ws.Range("E20").Value = '  +8.55%  '

# Row 21
$ws.Range("D21").Value = '''6.13'
$ws.Range("E21").Value = '  +6.64%  '

# Row 22
$ws.Range("D22").Value = '''66.79'
$ws.Range("E22").Value = '  +7.62%  '

# Row 23
$ws.Range("D23").Value = '''248.67'
$ws.Range("E23").Value = '  +9.73%  '

# Row 24
$ws.Range("D24").Value = '''2.91'
$ws.Range("E24").Value = '  +6.70%  '

# Row 25
$ws.Range("E25").Value = '  +10.91%  '

# Row 26
$ws.Range("E26").Value = '  -0.19%  '

# Row 27
$ws.Range("D27").Value = '''42.62'
$ws.Range("E27").Value = '  +23.42%  '

# Row 28
$ws.Range("E28").Value = '  +2.26%  '

# Row 29
$ws.Range("D29").Value = '''9.85'
$ws.Range("E29").Value = '  +8.80%  '

# Row 30
$ws.Range("D30").Value = '''20.20'
$ws.Range("E30").Value = '  +6.99%  '

# Row 31
$ws.Range("E31").Value = '  +10.99%  '

# Row 32
$ws.Range("D32").Value = '''147.17'
$ws.Range("E32").Value = '  +1.55%  '

# Row 33
$ws.Range("D33").Value = '''0.0798'
$ws.Range("E33").Value = '  +11.75%  '

# Row 34
$ws.Range("D34").Value = '''2.62'
$ws.Range("E34").Value = '  +6.20%  '

# Row 35
$ws.Range("E35").Value = '  +13.43%  '

# Row 36
$ws.Range("E36").Value = '  +9.57%  '

# Row 37
$ws.Range("E37").Value = '  +3.90%  '

# Row 38
$ws.Range("E38").Value = '  +11.51%  '

# Row 39
$ws.Range("D39").Value = '''15.75'
$ws.Range("E39").Value = '  +22.11%  '

# Row 40
$ws.Range("E40").Value = '  +16.50%  '

# Row 41
$ws.Range("D41").Value = '''3.40'
$ws.Range("E41").Value = '  +14.04%  '

# Row 42
$ws.Range("D42").Value = '''0.0307'
$ws.Range("E42").Value = '  +11.69%  '

# Row 43
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").Value = '''1.00'
$ws.Range("E43").Value = '  -0.51%  '

# Row 44
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").Value = '''1.98'
$ws.Range("E44").Value = '  +23.44%  '

# Row 45
$ws.Range("D45").Value = '1.845.54'
$ws.Range("E45").Value = '  +5.54%  '

# Row 46
$ws.Range("D46").Value = '''90.86'
$ws.Range("E46").Value = '  +26.52%  '

# Row 47
$ws.Range("D47").Value = '''0.200'
$ws.Range("E47").Value = '  +19.03%  '

# Row 48
$ws.Range("D48").Value = '''76.07'
$ws.Range("E48").Value = '  +19.31%  '

# Row 49
$ws.Range("D49").Value = '''4.96'
$ws.Range("E49").Value = '  +12.19%  '

# Row 50
$ws.Range("E50").Value = '  +7.81%  '

# Row 51
$ws.Range("D51").Value = '''54.22'
$ws.Range("E51").Value = '  +13.65%  '
